# Trade #92 closed at 2026-02-17 09:08:58 - unknown UNKNOWN +0.000%
#
# Appends the new trade row (trade #92 / row 93) to the "All Trades" and
# "MarketMaking" sheets, and rolls the new trade's results up into the
# "Summary" and "Strategy Status" sheets.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$status  = $wb.Worksheets.Item("Strategy Status")
$allTrades = $wb.Worksheets.Item("All Trades")
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# ---------------------------------------------------------------------
# New trade row data (row 93 on both "All Trades" and "MarketMaking")
# ---------------------------------------------------------------------
$newRow = 93
$tradeData = @{
    A = 92
    B = "2026-02-17"
    C = "09:08:52"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.15
    G = 0.13
    H = "CLOSED"
    I = -13.3333
    J = -0.02
    K = 99.64
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheet in @($allTrades, $marketMaking)) {
    $sheet.Range("A$newRow").Value = $tradeData.A

    # Force the Date column to text so Excel doesn't reinterpret the
    # "yyyy-mm-dd" string as a date serial number (it's stored as plain
    # text elsewhere in the column).
    $sheet.Range("B$newRow").NumberFormat = "@"
    $sheet.Range("B$newRow").Value = $tradeData.B

    $sheet.Range("C$newRow").Value = $tradeData.C
    $sheet.Range("D$newRow").Value = $tradeData.D
    $sheet.Range("E$newRow").Value = $tradeData.E
    $sheet.Range("F$newRow").Value = $tradeData.F
    $sheet.Range("G$newRow").Value = $tradeData.G
    $sheet.Range("H$newRow").Value = $tradeData.H
    $sheet.Range("I$newRow").Value = $tradeData.I
    $sheet.Range("J$newRow").Value = $tradeData.J
    $sheet.Range("K$newRow").Value = $tradeData.K
    $sheet.Range("L$newRow").Value = $tradeData.L
    $sheet.Range("M$newRow").Value = $tradeData.M
    $sheet.Range("N$newRow").Value = $tradeData.N
    $sheet.Range("O$newRow").Value = $tradeData.O
    $sheet.Range("P$newRow").Value = $tradeData.P
    $sheet.Range("Q$newRow").Value = $tradeData.Q
}

# ---------------------------------------------------------------------
# Roll the new trade up into "Summary"
# ---------------------------------------------------------------------
$summary.Range("B3").Value = 1199.64   # Current Capital
$summary.Range("B4").Value = -0.35     # Total P&L $
$summary.Range("B5").Value = -0.08     # Total P&L %
$summary.Range("B6").Value = 92        # Total Trades
$summary.Range("B8").Value = 38        # Losing Trades
$summary.Range("B9").Value = 41.3      # Win Rate %

# ---------------------------------------------------------------------
# Roll the new trade up into "Strategy Status" (MarketMaking row, r4)
# ---------------------------------------------------------------------
$status.Range("C4").Value = 99.64      # Capital
$status.Range("D4").Value = 92         # Trades
$status.Range("E4").Value = -0.35      # P&L $
$status.Range("F4").Value = -0.36      # P&L %
$status.Range("G4").Value = 41.3       # Win Rate %
